# Weekly price-sheet update: insert a new record as row 282 (pushing the
# existing rows 282-299 down to 283-300, so the sheet grows from 299 to
# 300 data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 282; this shifts rows
# 282..299 down to 283..300 and extends the used range / dimension.
$ws.Rows.Item(282).Insert()

# Populate the newly inserted row 282 with the new week's record.
$ws.Range("A282").Value = 10
$ws.Range("B282").Value = 'Vega Modelo de Temuco'
$ws.Range("C282").Value = 'La Araucanía'
$ws.Range("D282").Value = 44516
$ws.Range("E282").Value = 9
$ws.Range("F282").Value = 100112032
$ws.Range("G282").Value = 'Zapallo italiano'
$ws.Range("H282").Value = 'Sin especificar'
$ws.Range("I282").Value = 'Primera'
$ws.Range("J282").Value = 65
$ws.Range("K282").Value = 9000
$ws.Range("L282").Value = 9000
$ws.Range("M282").Value = 9000
$ws.Range("N282").Value = '$/caja 60 unidades'
$ws.Range("O282").Value = 'Región de Arica y Parinacota'
$ws.Range("P282").Value = 150
$ws.Range("Q282").Value = 60
$ws.Range("R282").Value = 'Hortaliza'
